$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.688.90"
$ws.Range("E2").Value = "  -1.38%  "
$ws.Range("D3").Value = "1.594.94"
$ws.Range("E3").Value = "  -1.62%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.30"
$ws.Range("E5").Value = "  -1.23%  "
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0618"
$ws.Range("E8").Value = "  -1.54%  "
$ws.Range("E9").Value = "  -1.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.67"
$ws.Range("E10").Value = "  -0.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0836"
$ws.Range("E11").Value = "  -0.93%  "
$ws.Range("D12").Value = "1.818.65"
$ws.Range("E12").Value = "  -1.68%  "
$ws.Range("D13").Value = "1.599.82"
$ws.Range("E13").Value = "  -1.33%  "
$ws.Range("E14").Value = "  -2.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.524"
$ws.Range("E15").Value = "  -2.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.74"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("D17").Value = "26.658.08"
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "208.56"
$ws.Range("E19").Value = "  -2.33%  "
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.70"
$ws.Range("E21").Value = "  -1.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.25"
$ws.Range("E22").Value = "  -1.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.34"
$ws.Range("E23").Value = "  +0.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.90"
$ws.Range("E24").Value = "  -1.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.55"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -2.63%  "
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("E29").Value = "  -0.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0503"
$ws.Range("E30").Value = "  -1.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.16"
$ws.Range("E31").Value = "  -1.19%  "
$ws.Range("E32").Value = "  -2.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.669"
$ws.Range("E33").Value = "  -5.12%  "
$ws.Range("E34").Value = "  -2.16%  "
$ws.Range("D35").Value = "1.298.88"
$ws.Range("E35").Value = "  -3.85%  "
$ws.Range("E36").Value = "  -0.61%  "
$ws.Range("E37").Value = "  -4.71%  "
$ws.Range("E38").Value = "  -2.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.836"
$ws.Range("E39").Value = "  -0.41%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.794"
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.20"
$ws.Range("E42").Value = "  -1.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.36"
$ws.Range("E43").Value = "  +0.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.58"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").Value = "1.730.98"
$ws.Range("E45").Value = "  -1.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.903"
$ws.Range("E46").Value = "  +6.01%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.65"
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "89.87"
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0984"
$ws.Range("E49").Value = "  -1.70%  "
$ws.Range("E50").Value = "  -1.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.52"
$ws.Range("E51").Value = "  -1.36%  "
